$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "sure bitti"
$ws.Range("E4").Select()
